$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: add two new columns and copy header formatting
$ws.Range("C1").Value = 'Unnamed: 2'
$ws.Range("D1").Value = 'Unnamed: 3'
$ws.Range("A1:B1").Copy()
$ws.Range("C1:D1").PasteSpecial(-4122)

# Data rows: cluster name / active cases pairs
$data = @(
    @('21 Moore Street Apartment Complex Moonee Ponds', 5),
    @('574 Plummer Street Building B1 & B2 Port Melbourne', 10),
    @('A Treat of France Cafe Carlton North', 5),
    @('Acquire BPO Southbank', 38),
    @('Al Haj Halal Meats Glenroy', 71),
    @('Al-Taqwa College Truganina', 30),
    @('Amiga Montessori Craigieburn', 19),
    @('CMV Truck and Bus Dealership Springvale', 5),
    @('CS Square Caroline Springs', 8),
    @('Campbellfield Heights Primary School Campbellfield', 6),
    @('Can Panel Cambellfield', 6),
    @('Cedars Medical Clinic Coburg', 39),
    @('Chemist Warehouse Campbellfield DC', 5),
    @('City of Greater Shepparton Community', 8),
    @('City of Hobsons Bay Community', 10),
    @('City of Moreland Social Network', 5),
    @('City of Wyndham Community', 8),
    @('Classy Cabinets and Kitchens Craigieburn', 18),
    @('Coles Barkly Square Brunswick August', 5),
    @('Coles Broadmeadows Central Shopping Centre', 9),
    @('Coles Coburg North Village Aug', 7),
    @('Coles Greenvale Shopping Centre', 5),
    @('Coles Roxburgh Village Roxburgh Park', 17),
    @('Community Kids Meadow Heights', 22),
    @('Concept Caravans Campbellfield', 5),
    @('Costco Wholesale Epping', 18),
    @('Crossroads Logistics Sunshine North', 7),
    @('Croydon Orthodontics', 8),
    @('Ernst and Young 8 Exhibition Street Melbourne', 12),
    @('Fitzroy Community School Fitzroy North', 48),
    @('Glenroy West Primary School', 7),
    @('Gloria Jeans Coffees Broadmeadows Central', 5),
    @('Health Care Providers Association South Melbourne', 12),
    @('IGA Meadow Heights Shopping Centre Meadow Heights', 5),
    @('Ilim College Glenroy Campus Hadfield', 25),
    @('Ilim Learning Sanctuary Glenroy', 14),
    @('Industrial Galvanizers Valmont Coatings Campbellfield', 13),
    @('Islamic College of Melbourne Tarneit', 6),
    @('KFC Fawkner', 21),
    @('Kasr Sweets Coolaroo', 25),
    @('Kids House Early Learning Cheltenham', 6),
    @('Learning Nest Early Learning Centre Meadow Heights', 9),
    @('Malvern Health and Fitness Clinic Malvern', 7),
    @('Metricon Homes Lot 1815 Canning Dr Mickleham', 5),
    @('Montessori Beginnings Greenvale', 6),
    @('MyCentre Childcare Broadmeadows', 25),
    @('Newbury Child and Community Centre Craigieburn', 10),
    @('Newport Football Club Altona North', 7),
    @('Nido Early School Moonee Ponds', 18),
    @('Nino Early Learning Adventures Lalor', 5),
    @('Northern Health Northern Hospital Epping Emergency Department Tier 1B', 39),
    @('OnQ Plumbing and Excavations Cragieburn', 7),
    @('Oporto Coolaroo', 13),
    @('Paisley Park Early Learning Centre Bundoora', 7),
    @('Panorama Construction Site Whitehorse Rd Box Hill', 76),
    @('Pascoe Vale Primary School Pascoe Vale', 7),
    @('Pearl Street Child Care Centre Glenroy', 11),
    @('People First Healthcare Home Residence Disability Support Taylors Lakes', 5),
    @('Serco Mill Park', 9),
    @('Southern Cross Station Crew Room Tissue Box Docklands', 5),
    @('St Vincents Hospital Emergency Department Melbourne', 6),
    @('Tek Foods Somerton', 7),
    @('The Homestead Child and Family Centre Roxburgh Park', 17),
    @('The Royal Children''s Hospital Melbourne Emergency Department Parkville Tier 1B', 5),
    @('Tip Top Warehouse Dandenong', 9),
    @('Total Window Concepts Hoppers Crossing', 7),
    @('Turnstall Fresh Turnstall Square Shopping Centre Donaster East', 5),
    @('Victorian Civil and Administrative Tribunal Melbourne', 5),
    @('Werribee Mercy Hospital Emergency Department', 7),
    @('Western Health Footscray Hospital Emergency Department', 7),
    @('Woolworths Coburg Station', 5),
    @('Woolworths Greenvale Lakes Roxburgh Park', 15),
    @('Yara Childcare Centre Truganina', 31)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}

# Remove the now-unused trailing rows (sheet shrank from 76 to 74 rows)
$ws.Rows.Item(75).Delete()
$ws.Rows.Item(75).Delete()
